# Rename header columns from "Soybean.*" to "Mungbean.*" to reflect that
# this observed-data workbook now documents a Mungbean model (not soybean).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observed")

$ws.Range("C1").Value = "Mungbean.Phenology.CurrentStageName"
$ws.Range("D1").Value = "Mungbean.Phenology.MaturityDAS"
$ws.Range("E1").Value = "Mungbean.AboveGround.Wt"
$ws.Range("F1").Value = "Mungbean.AboveGround.Wterror"
$ws.Range("G1").Value = "Mungbean.Grain.Wt"
$ws.Range("H1").Value = "Mungbean.Grain.Wterror"
$ws.Range("K1").Value = "Mungbean.Grain.HarvestIndex"
$ws.Range("L1").Value = "Mungbean.Grain.HarvestIndexerror"

# Match the cursor/selection position recorded in the saved file.
$ws.Range("K18").Select()
